$d = $word.ActiveDocument

# =======================================================================
# Edit 1: Add a "Kommentar zu Abbildung f)" heading plus an explanatory
# paragraph right before the "Pointer-Semantik" section.
#
# In the original document there is an *empty* Heading-1 paragraph
# immediately followed by a page-break paragraph and then the
# "Pointer-Semantik" heading. We locate that empty heading paragraph
# structurally (empty Heading-1 paragraph, two paragraphs before a
# paragraph that starts with "Pointer-Semantik") instead of relying on
# a fixed paragraph index.
# =======================================================================

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text -eq [char]13) {
        if (($i + 2) -le $d.Paragraphs.Count) {
            $next2 = $d.Paragraphs.Item($i + 2)
            if ($next2.Range.Text.StartsWith("Pointer-Semantik")) {
                $targetIndex = $i
                break
            }
        }
    }
}
if ($targetIndex -eq -1) {
    throw "edit.ps1: could not locate the empty heading paragraph before 'Pointer-Semantik'"
}

$emptyHeadingPara = $d.Paragraphs.Item($targetIndex)
$emptyHeadingPara.Range.Text = "Kommentar zu Abbildung f)"

# Insert a brand-new (default/"Standard"-styled) paragraph right after the
# heading, containing the explanatory text about the overflow observed in
# figure (f).
$emptyHeadingPara.Range.InsertParagraphAfter()
$explPara = $d.Paragraphs.Item($targetIndex + 1)
$explPara.Range.Style = "Standard"

$er = $explPara.Range
$er.Collapse(1)
$er.InsertAfter("In der Abbildung (f) ist ein Overflow des Ringbuffers zu beobachten. Es wird also ein Element mehr rein geschrieben, als eigentlich vom Ringbuffer gehalten werden kann. Somit wird das älteste Element des Ringbuffer mit dem neuen Element")
$er = $d.Range($er.End, $er.End)
$er.InsertAfter(" ")
$er = $d.Range($er.End, $er.End)
$er.InsertAfter("überschrieben. Sowohl In- als auch Out-Pointer werden dabei einen weiter gesetzt, um auf das nun älteste Element zu zeigen, um den „normalen“ Ablauf weiter zu betreiben.")

# =======================================================================
# Edit 2: In the "Pointer-Semantik" explanation paragraph, add two new
# sentences about write operations on the ring buffer (right after
# "... Speicherzelle gesetzt. " and before "Ist einer der beiden
# Pointer ..."), with "RB", "p_out" and "p_in" set in Courier New like
# the rest of the code-styled terms in that paragraph.
# =======================================================================

$anchor = $d.Content
$anchor.Find.ClearFormatting()
$ok = $anchor.Find.Execute("Ist einer der beiden Pointer im Begriff", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "edit.ps1: could not locate edit-2 anchor text 'Ist einer der beiden Pointer im Begriff'"
}
$insPos = $anchor.Start

$fullText2 = "Wird eine Schreiboperation auf RB ausgeführt, wird im Anschluss p_out auf die Adresse der nachfolgenden Speicherzelle gesetzt. Außerdem darf p_in nicht p_out überholen, da sonst nicht gelesene Zeichen überschrieben würden. "

$insRange = $d.Range($insPos, $insPos)
$insRange.InsertAfter($fullText2)

function SetMonoAt($baseOffset, $text, $needle) {
    $idx = $text.IndexOf($needle)
    while ($idx -ge 0) {
        $s = $baseOffset + $idx
        $e = $s + $needle.Length
        $mr = $d.Range($s, $e)
        $mr.Font.Name = "Courier New"
        $mr.Font.Size = 10
        $idx = $text.IndexOf($needle, $idx + $needle.Length)
    }
}

SetMonoAt $insPos $fullText2 "RB"
SetMonoAt $insPos $fullText2 "p_out"
SetMonoAt $insPos $fullText2 "p_in"

Write-Output "edit.ps1 completed successfully"
